$d = $word.ActiveDocument

# Update the date heading at the top of the document
$d.Content.Find.Execute("2024-04-29 Monday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-04-30 Tuesday", 2)

# Update the division problems in the table, cell by cell (positional,
# since several values repeat and a blind Find/Replace would be ambiguous).
$t = $d.Tables.Item(1)

$newValues = @{
    1  = @("23÷4=", "60÷2=", "85÷3=", "77÷7=", "34÷9=")
    5  = @("79÷2=", "35÷5=", "85÷8=", "23÷4=", "69÷7=")
    9  = @("85÷3=", "68÷4=", "41÷3=", "62÷4=", "33÷4=")
    13 = @("86÷9=", "27÷6=", "27÷8=", "61÷5=", "32÷9=")
    17 = @("48÷6=", "22÷5=", "24÷9=", "73÷4=", "94÷5=")
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($col = 1; $col -le $vals.Length; $col++) {
        $cell = $t.Cell($row, $col)
        $cell.Range.Text = $vals[$col - 1]
    }
}
